$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPR")

# Insert two new columns before column D; existing D:K data shifts to F:M.
$ws.Columns("D:E").Insert()

# The newly inserted D:E columns should carry the same formatting as the rest
# of each row (date format on header rows, number format on data rows); copy
# that formatting in bulk from column F (which now holds the old column D values).
# Work block-by-block so we never touch rows that have no D:K cells at all
# (section-title rows 5, 6, 37, 79 and the blank separator rows 36, 78), which
# must stay exactly as they were.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("F7:F35").Copy()
$ws.Range("E7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns with the latest two quarters of data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 130700
$ws.Range("E8").Value = 131100
$ws.Range("D9").Value = 19000
$ws.Range("E9").Value = 20100
$ws.Range("D10").Value = 111700
$ws.Range("E10").Value = 111000
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 6500
$ws.Range("E14").Value = 5100
$ws.Range("D15").Value = 76400
$ws.Range("E15").Value = 58900
$ws.Range("D17").Value = 114500
$ws.Range("E17").Value = 96200
$ws.Range("D18").Value = 16200
$ws.Range("E18").Value = 34900
$ws.Range("D20").Value = 221400
$ws.Range("E20").Value = -51100
$ws.Range("D21").Value = 314000
$ws.Range("E21").Value = 42800
$ws.Range("D22").Value = 13400
$ws.Range("E22").Value = 13200
$ws.Range("D23").Value = 224300
$ws.Range("E23").Value = -29400
$ws.Range("D24").Value = 1800
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 222400
$ws.Range("E26").Value = -29400
$ws.Range("D27").Value = 222400
$ws.Range("E27").Value = -29400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -221400
$ws.Range("E32").Value = 51100
$ws.Range("D33").Value = 222400
$ws.Range("E33").Value = -29400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 222400
$ws.Range("E35").Value = -29400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 32800
$ws.Range("E41").Value = 93000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 72900
$ws.Range("E43").Value = 68000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 84100
$ws.Range("E45").Value = 3000
$ws.Range("D46").Value = 189800
$ws.Range("E46").Value = 164000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2029500
$ws.Range("E48").Value = 1973900
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 33200
$ws.Range("E52").Value = 6800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2252500
$ws.Range("E54").Value = 2144600
$ws.Range("D57").Value = 86800
$ws.Range("E57").Value = 66800
$ws.Range("D58").Value = 1900
$ws.Range("E58").Value = 2000
$ws.Range("D59").Value = 159500
$ws.Range("E59").Value = 270000
$ws.Range("D60").Value = 248200
$ws.Range("E60").Value = 338800
$ws.Range("D61").Value = 617400
$ws.Range("E61").Value = 617000
$ws.Range("D62").Value = 174800
$ws.Range("E62").Value = 201000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1040400
$ws.Range("E66").Value = 1156800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -559800
$ws.Range("E72").Value = -782300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1212100
$ws.Range("E76").Value = 987800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 222400
$ws.Range("E81").Value = -29400
$ws.Range("D83").Value = 76400
$ws.Range("E83").Value = 58900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 71300
$ws.Range("E89").Value = 91300
$ws.Range("D91").Value = -200
$ws.Range("E91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -131100
$ws.Range("E94").Value = -102500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -400
$ws.Range("E100").Value = -3200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -60200
$ws.Range("E102").Value = -14400

# Row 91 (Capital Expenditures) also received corrected figures for the older
# quarters F:J, not just the two newly inserted columns.
$ws.Range("F91").Value = -300
$ws.Range("G91").Value = -100
$ws.Range("H91").Value = -700
$ws.Range("I91").Value = -100
$ws.Range("J91").Value = -200
